# Update "想去人数" (interest count) values in the "展览" sheet and the
# aggregated "全部类型" sheet to reflect the latest scrape at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Worksheet "展览" (rows 2-10, column F)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 5508
$wsExhibit.Range("F3").Value = 602
$wsExhibit.Range("F4").Value = 12174
$wsExhibit.Range("F5").Value = 300
$wsExhibit.Range("F6").Value = 612
$wsExhibit.Range("F7").Value = 182
$wsExhibit.Range("F8").Value = 337
$wsExhibit.Range("F9").Value = 1113
$wsExhibit.Range("F10").Value = 105

# Worksheet "全部类型" (same events, different rows, column F)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 5508
$wsAll.Range("F4").Value = 602
$wsAll.Range("F6").Value = 12174
$wsAll.Range("F7").Value = 300
$wsAll.Range("F8").Value = 612
$wsAll.Range("F9").Value = 182
$wsAll.Range("F12").Value = 337
$wsAll.Range("F13").Value = 1113
$wsAll.Range("F15").Value = 105
